$d = $word.ActiveDocument

# --- 1) Title-page date: "January 11, 2010" -> "January 11" + superscript "th" + ", 2010" ---

# Replace the plain date text with date+superscript-suffix placeholder text first.
# Replace:=1 (wdReplaceOne) so only the FIRST hit (the title-page date) is touched -
# the Revision History table further down also contains "January 11, 2010" and must
# stay untouched.
$null = $d.Content.Find.Execute("January 11, 2010", $true, $false, $false, $false, $false, $true, 1, $false, "January 11th, 2010", 1)

# Now locate the "th" that immediately follows "January 11" so we only touch that occurrence
# (the Revision History table also contains "January 11, 2010" but is untouched because it
# never got the "th" inserted into it).
$thRng = $d.Content
$null = $thRng.Find.Execute("th, 2010")
$thOnly = $d.Range($thRng.Start, $thRng.Start + 2)
$thOnly.Font.Superscript = $true

# --- 2) Footer page-number field cached text: "1" -> "3" (body section footer) ---

$bodyFooter = $d.Sections.Item(2).Footers.Item(1)
$null = $bodyFooter.Range.Find.Execute("1", $false, $false, $false, $false, $false, $true, 1, $false, "3", 1)
